$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 ("Marking"): correct the per-question marking values
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): correct the total marks after fixing the marking scheme
$ws.Range("B12").Value = 108
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "106 / 112"
